$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row so each column carries the respective
#        input-file-name suffix: "_old" -> "_FV2310", "_new" -> "_FV2404".
#        The "diff" header (column K) is left untouched.
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2. Turn the data range into a real Excel Table ("Table1") so the
#        header row gets the filter buttons / structured reference support.
$tblRange = $ws.Range("A1:U87")
$tbl = $ws.ListObjects.Add(1, $tblRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1, top-left of the
#        scrollable pane is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
